$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 3) mirroring the structure of row 2, with new values.
$ws.Range("A3").Value = 33
$ws.Range("B3").Value = "test archivage"
$ws.Range("C3").Value = "2020-02-28 00:00:00"
$ws.Range("D3").Value = "2020-03-01 12:10:00"
$ws.Range("E3").Value = "2020-03-11 12:10:00"
$ws.Range("F3").Value = "Test d'archivage"
$ws.Range("G3").Value = "Drakyn"
$ws.Range("H3").Value = "Saint Herblain"
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = "Place de la fontaine Lille (59000)"
